$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Create_EPG example row) - rename from Storage_Mgmt to Py test values
$ws.Range("B2").Value = "PythonTest"
$ws.Range("C2").Value = "Py_test1"
$ws.Range("D2").Value = "Py_mgmt"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "PythonTest_BD"

# Row 3 - clear leftover "Create_EPG" label in column A (row 3 no longer a separate example)
$ws.Range("A3").Value = ""

# Row 6 (Create_BD example row) - rename from Storage_Mgmt_BD / Prod to Py test values
$ws.Range("B6").Value = "PythonTest"
$ws.Range("C6").Value = "PythonTest_BD"
$ws.Range("D6").Value = "Py_Prod"
$ws.Range("E6").Value = "10.207.250.1/24"
$ws.Range("F6").Value = "no"

# Update selected cell to reflect the saved view state
$ws.Range("H2").Select()
